$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $newValue) {
    $rng = $ws.Range($cellRef)
    $origStyle = $rng.Style
    $rng.NumberFormat = '@'
    $rng.Value = $newValue
    $rng.Style = $origStyle
}

Set-TextValue 'D2' '23.441.50'
Set-TextValue 'E2' '  -1.04%  '
Set-TextValue 'D3' '1.646.26'
Set-TextValue 'E3' '  -0.56%  '
Set-TextValue 'D4' '0.9982'
Set-TextValue 'E4' '  -0.41%  '
Set-TextValue 'D5' '0.9985'
Set-TextValue 'E5' '  -0.31%  '
Set-TextValue 'D6' '299.61'
Set-TextValue 'E6' '  -1.08%  '
Set-TextValue 'D7' '0.3795'
Set-TextValue 'E7' '  -1.02%  '
Set-TextValue 'D8' '50.45'
Set-TextValue 'E8' '  -1.29%  '
Set-TextValue 'D9' '0.3497'
Set-TextValue 'E9' '  -3.01%  '
Set-TextValue 'D10' '0.08076'
Set-TextValue 'E10' '  -1.44%  '
Set-TextValue 'D11' '1.220'
Set-TextValue 'E11' '  -0.87%  '
Set-TextValue 'D12' '0.9980'
Set-TextValue 'E12' '  -0.41%  '
Set-TextValue 'D13' '22.09'
Set-TextValue 'E13' '  -1.41%  '
Set-TextValue 'D14' '6.314'
Set-TextValue 'E14' '  -2.18%  '
Set-TextValue 'D15' '7.279'
Set-TextValue 'E15' '  -2.23%  '
Set-TextValue 'D16' '0.00001218'
Set-TextValue 'E16' '  -0.32%  '
Set-TextValue 'D17' '1.640.01'
Set-TextValue 'E17' '  -0.89%  '
Set-TextValue 'D18' '95.00'
Set-TextValue 'E18' '  -2.66%  '
Set-TextValue 'D19' '0.06966'
Set-TextValue 'E19' '  -1.02%  '
Set-TextValue 'D20' '6.633'
Set-TextValue 'E20' '  -2.21%  '
Set-TextValue 'D21' '17.37'
Set-TextValue 'E21' '  -1.07%  '
Set-TextValue 'D22' '0.9982'
Set-TextValue 'E22' '  -0.36%  '
Set-TextValue 'D23' '12.43'
Set-TextValue 'E23' '  -2.16%  '
Set-TextValue 'D24' '23.442.76'
Set-TextValue 'E24' '  -1.02%  '
Set-TextValue 'D25' '2.437'
Set-TextValue 'E25' '  -1.99%  '
Set-TextValue 'D26' '2.983'
Set-TextValue 'E26' '  -1.45%  '
Set-TextValue 'D27' '21.02'
Set-TextValue 'E27' '  -1.13%  '
Set-TextValue 'D28' '149.91'
Set-TextValue 'E28' '  -1.93%  '
Set-TextValue 'D29' '5.179'
Set-TextValue 'E29' '  -1.14%  '
Set-TextValue 'D30' '131.60'
Set-TextValue 'D31' '1.819.76'
Set-TextValue 'E31' '  -1.14%  '
Set-TextValue 'D32' '6.825'
Set-TextValue 'E32' '  -3.74%  '
Set-TextValue 'D33' '2.132'
Set-TextValue 'E33' '  -5.42%  '
Set-TextValue 'D34' '11.28'
Set-TextValue 'E34' '  -5.92%  '
Set-TextValue 'D35' '0.9887'
Set-TextValue 'E35' '  -6.53%  '
Set-TextValue 'E36' '  -4.28%  '
Set-TextValue 'D37' '0.08779'
Set-TextValue 'E37' '  -0.33%  '
Set-TextValue 'D38' '0.2426'
Set-TextValue 'E38' '  -3.06%  '
Set-TextValue 'D39' '5.888'
Set-TextValue 'E39' '  -3.16%  '
Set-TextValue 'D40' '0.06831'
Set-TextValue 'E40' '  -2.30%  '
Set-TextValue 'E41' '  -1.89%  '
Set-TextValue 'D42' '0.6839'
Set-TextValue 'E42' '  -2.13%  '
Set-TextValue 'D43' '1.292'
Set-TextValue 'E43' '  -3.43%  '
Set-TextValue 'D44' '15.61'
Set-TextValue 'E44' '  -2.00%  '
Set-TextValue 'D45' '0.9974'
Set-TextValue 'E45' '  -0.38%  '
Set-TextValue 'D46' '0.6348'
Set-TextValue 'E46' '  -2.42%  '
Set-TextValue 'D47' '2.242'
Set-TextValue 'E47' '  -2.44%  '
Set-TextValue 'D48' '3.912'
Set-TextValue 'E48' '  -1.31%  '
Set-TextValue 'B49' 'Cronos'
Set-TextValue 'C49' 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
Set-TextValue 'D49' '0.07684'
Set-TextValue 'E49' '  -2.72%  '
Set-TextValue 'B50' 'Quant'
Set-TextValue 'C50' 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
Set-TextValue 'D50' '127.18'
Set-TextValue 'E50' '  -0.83%  '
Set-TextValue 'D51' '1.219'
Set-TextValue 'E51' '  +2.09%  '
